$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.864.21'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.218.16'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.02'
$ws.Range("E5").Value = '  +4.91%  '
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '77.44'
$ws.Range("E7").Value = '  +2.86%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("E9").Value = '  -1.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.61'
$ws.Range("E10").Value = '  +3.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0909'
$ws.Range("E11").Value = '  -2.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.07'
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.551.42'
$ws.Range("E14").Value = '  -1.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.47'
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.208.98'
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.785'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.865.02'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("E19").Value = '  -1.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.19'
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.99'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.25'
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.72'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.37'
$ws.Range("E24").Value = '  -5.35%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '42.80'
$ws.Range("E26").Value = '  +10.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.80'
$ws.Range("E28").Value = '  -1.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.00'
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.42'
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0876'
$ws.Range("E33").Value = '  +10.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.23'
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0359'
$ws.Range("E36").Value = '  +7.95%  '
$ws.Range("E37").Value = '  -1.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.34'
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.04'
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("E40").Value = '  +17.64%  '
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.200'
$ws.Range("E42").Value = '  -2.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.32'
$ws.Range("E43").Value = '  -2.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.48'
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.33'
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("B46").Value = 'WOONetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.472'
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.38'
$ws.Range("E47").Value = '  -3.71%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0977'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.437.48'
$ws.Range("E51").Value = '  -0.81%  '
